# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-04 10:15:24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage figures in this report are stored as literal text (e.g. "67.0%"),
# not real numeric percentages, so force a Text number format before writing
# them - otherwise Excel's COM layer helpfully (and wrongly, for this sheet)
# reinterprets the string as a numeric percentage value.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# Workbook-level summary metrics (K/L columns, rows 4-10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 213          # Recorded Sessions
$ws.Range("L7").Value = 9            # Missing Sessions
Set-TextValue $ws.Range("L9") "67.0%"    # Coverage %
Set-TextValue $ws.Range("L10") "76.0%"   # Average Attendance %

# ---------------------------------------------------------------------------
# "Recorded By" column (G) swaps from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" for every session row that had both recorders.
# ---------------------------------------------------------------------------
$gSwapRows = @(8, 9, 10, 12, 14, 15, 17, 18, 34, 35, 36, 38, 40, 41, 43, 44, 60, 61, 62, 64, 66, 67, 69, 70, 86, 87, 88, 90, 92, 93, 95, 96, 112, 113, 114, 116, 118, 119, 121, 122, 138, 139, 140, 142, 144, 145, 147, 148, 164, 167, 170, 174, 191, 194, 197, 201, 218, 221, 224, 228, 245, 248, 251, 255, 272, 275, 278, 282, 299, 302, 305, 309)
foreach ($r in $gSwapRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# Group Statistics block (K14:S20) - recorded/missing counts and percentages
# shifted by newly-recorded sessions.
# ---------------------------------------------------------------------------
# Row 15 - B1A1
$ws.Range("O15").Value = 18
$ws.Range("P15").Value = 1
Set-TextValue $ws.Range("R15") "69.2%"
Set-TextValue $ws.Range("S15") "81.6%"

# Row 16 - B1A2
$ws.Range("O16").Value = 19
$ws.Range("P16").Value = 0
Set-TextValue $ws.Range("R16") "73.1%"
Set-TextValue $ws.Range("S16") "80.7%"

# Row 17 - B1B1
$ws.Range("O17").Value = 19
$ws.Range("P17").Value = 0
Set-TextValue $ws.Range("R17") "73.1%"

# Row 18 - B1B2
$ws.Range("O18").Value = 19
$ws.Range("P18").Value = 0
Set-TextValue $ws.Range("R18") "73.1%"
Set-TextValue $ws.Range("S18") "78.2%"

# Row 19 - B1C1
$ws.Range("O19").Value = 19
$ws.Range("P19").Value = 0
Set-TextValue $ws.Range("R19") "73.1%"
Set-TextValue $ws.Range("S19") "76.8%"

# Row 20 - B1C2
$ws.Range("O20").Value = 18
$ws.Range("P20").Value = 1
Set-TextValue $ws.Range("R20") "69.2%"
Set-TextValue $ws.Range("S20") "80.0%"

# ---------------------------------------------------------------------------
# Session rows that moved from "Not Recorded" to "Recorded" (row 20 of each
# group block, 20 = B1A1, 46 = B1A2, 72 = B1B1, 98 = B1B2, 124 = B1C1,
# 150 = B1C2). Each needs: recorder email filled in, the attended/total
# students count updated, status text, and the pink "Not Recorded" fill
# swapped for the green "Recorded" fill used elsewhere in the sheet. The
# green fill/font/alignment formatting is pulled straight from the
# already-"Recorded" row directly above via PasteSpecial (formats only),
# which re-uses the workbook's existing style instead of fabricating a
# new one.
# ---------------------------------------------------------------------------
$recordedRows = @(
    @{ Row = 20;  FormatSourceRow = 19;  Students = "20/26" },
    @{ Row = 46;  FormatSourceRow = 45;  Students = "20/27" },
    @{ Row = 72;  FormatSourceRow = 71;  Students = "18/26" },
    @{ Row = 98;  FormatSourceRow = 97;  Students = "27/27" },
    @{ Row = 124; FormatSourceRow = 123; Students = "25/30" },
    @{ Row = 150; FormatSourceRow = 149; Students = "18/23" }
)

foreach ($entry in $recordedRows) {
    $r = $entry.Row

    $src = $ws.Range("A" + $entry.FormatSourceRow + ":I" + $entry.FormatSourceRow)
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats

    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com"
    $ws.Cells.Item($r, 8).Value = $entry.Students
    $ws.Cells.Item($r, 9).Value = "Recorded"
}
